$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.427.25"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "1.851.66"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.22"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("E6").Value = "  -3.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07663"
$ws.Range("E8").Value = "  +2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2979"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.55"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").Value = "1.993.67"
$ws.Range("E11").Value = "  +7.46%  "

$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.010"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6901"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009998"
$ws.Range("E15").Value = "  +4.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.14"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "2.205.82"
$ws.Range("E17").Value = "  +4.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.202"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").Value = "29.609.55"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.86"
$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.674"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.91"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1392"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.486"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.477"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05787"
$ws.Range("E30").Value = "  -4.16%  "

$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.135"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.882"

$ws.Range("E35").Value = "  -1.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7212"
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").Value = "1.253.36"
$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.806"
$ws.Range("E39").Value = "  +0.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01809"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9103"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.095"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").Value = "2.138.49"
$ws.Range("E43").Value = "  +5.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.75"
$ws.Range("E45").Value = "  +1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.75"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.298"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000118"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.181"
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4038"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.702"
$ws.Range("E51").Value = "  +2.21%  "
